# Onmyoji Modifiers - update area-targeting selections on Sheet2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update the dropdown selections in row 3 (B3, C3, D3, F3). E3 stays "Roll Count".
$ws.Range("B3").Value = "All"
$ws.Range("C3").Value = "AoE"
$ws.Range("D3").Value = "All"
$ws.Range("F3").Value = "All"

# Move the active selection to I4, matching the reviewed output cell.
$ws.Activate()
$ws.Range("I4").Select()

$wb.Application.Calculate()
